$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value2 = 10108454
$ws.Range("I76").Value2 = 8037.65
$ws.Range("J76").Value2 = 25647556
$ws.Range("K76").Value2 = 8037.65
$ws.Range("L76").Value2 = 25647556
$ws.Range("M76").Value2 = -7722.65
$ws.Range("N76").Value2 = -25648186

$ws.Range("H79").Value2 = 10108454
$ws.Range("I79").Value2 = 8037.65
$ws.Range("J79").Value2 = 25647556
$ws.Range("K79").Value2 = 8037.65
$ws.Range("L79").Value2 = 25647556
$ws.Range("M79").Value2 = -6945.65
$ws.Range("N79").Value2 = -25649740

$ws.Range("H141").Value2 = 8844.429
$ws.Range("I141").Value2 = 4567
$ws.Range("J141").Value2 = 12052.5
$ws.Range("K141").Value2 = 13701
$ws.Range("L141").Value2 = 36157.5
$ws.Range("M141").Value2 = -8521
$ws.Range("N141").Value2 = -46517.5


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value2 = 19794.834
$ws.Range("J80").Value2 = 19794.834
$ws.Range("L80").Value2 = 19794.834
$ws.Range("N80").Value2 = -21790.834

$ws.Range("H83").Value2 = 19794.834
$ws.Range("J83").Value2 = 19794.834
$ws.Range("L83").Value2 = 59384.50199999999
$ws.Range("N83").Value2 = -69368.50199999999

$ws.Range("H102").Value2 = 1576.125
$ws.Range("I102").Value2 = 1444.1428
$ws.Range("K102").Value2 = 1444.1428
$ws.Range("M102").Value2 = 177.8571999999999

$ws.Range("H132").Value2 = 2675399.2
$ws.Range("I132").Value2 = 1143.3077
$ws.Range("J132").Value2 = 6538213
$ws.Range("K132").Value2 = 3429.9231
$ws.Range("L132").Value2 = 19614639
$ws.Range("M132").Value2 = -899.9231
$ws.Range("N132").Value2 = -19619699


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 2385.5
$ws.Range("I16").Value2 = 1187.5
$ws.Range("J16").Value2 = 2984.5
$ws.Range("K16").Value2 = 1187.5
$ws.Range("L16").Value2 = 2984.5
$ws.Range("M16").Value2 = -900.5
$ws.Range("N16").Value2 = -3558.5

$ws.Range("H31").Value2 = 3580.9211
$ws.Range("I31").Value2 = 858.5769
$ws.Range("J31").Value2 = 9479.333000000001
$ws.Range("K31").Value2 = 858.5769
$ws.Range("L31").Value2 = 9479.333000000001
$ws.Range("M31").Value2 = -563.5769
$ws.Range("N31").Value2 = -10069.333

$ws.Range("H34").Value2 = 3580.9211
$ws.Range("I34").Value2 = 858.5769
$ws.Range("J34").Value2 = 9479.333000000001
$ws.Range("K34").Value2 = 858.5769
$ws.Range("L34").Value2 = 9479.333000000001
$ws.Range("M34").Value2 = -656.5769
$ws.Range("N34").Value2 = -9883.333000000001

$ws.Range("H62").Value2 = 4243.3335
$ws.Range("J62").Value2 = 5477.2
$ws.Range("L62").Value2 = 5477.2
$ws.Range("N62").Value2 = -6725.2

$ws.Range("H65").Value2 = 4243.3335
$ws.Range("J65").Value2 = 5477.2
$ws.Range("L65").Value2 = 27386
$ws.Range("N65").Value2 = -33626

$ws.Range("H113").Value2 = 2385.5
$ws.Range("I113").Value2 = 1187.5
$ws.Range("J113").Value2 = 2984.5
$ws.Range("K113").Value2 = 1187.5
$ws.Range("L113").Value2 = 2984.5
$ws.Range("M113").Value2 = 982.5
$ws.Range("N113").Value2 = -7324.5


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 29917608
$ws.Range("I5").Value2 = 41667030
$ws.Range("J5").Value2 = 9977.182000000001
$ws.Range("K5").Value2 = 125001090
$ws.Range("L5").Value2 = 29931.546
$ws.Range("M5").Value2 = -125000978
$ws.Range("N5").Value2 = -30155.546

$ws.Range("H68").Value2 = 8115.154
$ws.Range("I68").Value2 = 314
$ws.Range("J68").Value2 = 17216.5
$ws.Range("K68").Value2 = 942
$ws.Range("L68").Value2 = 51649.5
$ws.Range("M68").Value2 = -131
$ws.Range("N68").Value2 = -53271.5

$ws.Range("H71").Value2 = 8115.154
$ws.Range("I71").Value2 = 314
$ws.Range("J71").Value2 = 17216.5
$ws.Range("K71").Value2 = 2826
$ws.Range("L71").Value2 = 154948.5
$ws.Range("M71").Value2 = 1230
$ws.Range("N71").Value2 = -163060.5

$ws.Range("H122").Value2 = 12504780
$ws.Range("I122").Value2 = 69444810
$ws.Range("J122").Value2 = 5750.1953
$ws.Range("K122").Value2 = 625003290
$ws.Range("L122").Value2 = 51751.7577
$ws.Range("M122").Value2 = -625000840
$ws.Range("N122").Value2 = -56651.7577

$ws.Range("H131").Value2 = 855.99
$ws.Range("J131").Value2 = 864.6288500000001
$ws.Range("L131").Value2 = 2593.88655
$ws.Range("N131").Value2 = -12673.88655

$ws.Range("H132").Value2 = 5299.56
$ws.Range("I132").Value2 = 771.8570999999999
$ws.Range("J132").Value2 = 7060.3335
$ws.Range("K132").Value2 = 6946.7139
$ws.Range("L132").Value2 = 63543.0015
$ws.Range("M132").Value2 = -4416.7139
$ws.Range("N132").Value2 = -68603.0015

$ws.Range("H135").Value2 = 29917608
$ws.Range("I135").Value2 = 41667030
$ws.Range("J135").Value2 = 9977.182000000001
$ws.Range("K135").Value2 = 375003270
$ws.Range("L135").Value2 = 89794.63800000001
$ws.Range("M135").Value2 = -375000735
$ws.Range("N135").Value2 = -94864.63800000001


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value2 = 3658.3333
$ws.Range("I31").Value2 = 487.5
$ws.Range("J31").Value2 = 10000
$ws.Range("K31").Value2 = 487.5
$ws.Range("L31").Value2 = 10000
$ws.Range("M31").Value2 = -195.5
$ws.Range("N31").Value2 = -10584

$ws.Range("H37").Value2 = 3658.3333
$ws.Range("I37").Value2 = 487.5
$ws.Range("J37").Value2 = 10000
$ws.Range("K37").Value2 = 487.5
$ws.Range("L37").Value2 = 10000
$ws.Range("M37").Value2 = -210.5
$ws.Range("N37").Value2 = -10554

$ws.Range("H52").Value2 = 39998
$ws.Range("I52").Value2 = 0
$ws.Range("J52").Value2 = 39998
$ws.Range("K52").Value2 = 0
$ws.Range("L52").Value2 = 39998
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value2 = -40516


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value2 = 0
$ws.Range("J14").Value2 = 0
$ws.Range("L14").Value2 = 0
$ws.Range("N14").ClearContents()

$ws.Range("H20").Value2 = 7727.273
$ws.Range("J20").Value2 = 7727.273
$ws.Range("L20").Value2 = 7727.273
$ws.Range("N20").Value2 = -8179.273

$ws.Range("H46").Value2 = 4168575
$ws.Range("J46").Value2 = 3420.5
$ws.Range("L46").Value2 = 3420.5
$ws.Range("N46").Value2 = -3796.5

$ws.Range("H132").Value2 = 51959870
$ws.Range("I132").Value2 = 87913050
$ws.Range("J132").Value2 = 27500.334
$ws.Range("K132").Value2 = 263739150
$ws.Range("L132").Value2 = 82501.00199999999
$ws.Range("M132").Value2 = -263736620
$ws.Range("N132").Value2 = -87561.00199999999

$ws.Range("H136").Value2 = 165083300
$ws.Range("I136").Value2 = 210888210
$ws.Range("J136").Value2 = 125004000
$ws.Range("K136").Value2 = 632664630
$ws.Range("L136").Value2 = 375012000
$ws.Range("M136").Value2 = -632662080
$ws.Range("N136").Value2 = -375017100

